# Auto-generated edit script applying the Raiden_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 852.86957
$ws.Range("I15").Value = 852.86957
$ws.Range("K15").Value = 2558.60871
$ws.Range("M15").Value = -2389.60871

$ws.Range("H80").Value = 419.15384
$ws.Range("I80").Value = 662
$ws.Range("J80").Value = 311.22223
$ws.Range("K80").Value = 1986
$ws.Range("L80").Value = 933.66669
$ws.Range("M80").Value = -988
$ws.Range("N80").Value = -2929.66669

$ws.Range("H83").Value = 419.15384
$ws.Range("I83").Value = 662
$ws.Range("J83").Value = 311.22223
$ws.Range("K83").Value = 5958
$ws.Range("L83").Value = 2801.00007
$ws.Range("M83").Value = -966
$ws.Range("N83").Value = -12785.00007

$ws.Range("H86").Value = 3350.3044
$ws.Range("J86").Value = 3549.5
$ws.Range("L86").Value = 3549.5
$ws.Range("N86").Value = -5795.5

$ws.Range("H89").Value = 3350.3044
$ws.Range("J89").Value = 3549.5
$ws.Range("L89").Value = 17747.5
$ws.Range("N89").Value = -28979.5

$ws.Range("H138").Value = 2755.8594
$ws.Range("I138").Value = 5012.3335
$ws.Range("J138").Value = 2235.1345
$ws.Range("K138").Value = 15037.0005
$ws.Range("L138").Value = 6705.4035
$ws.Range("M138").Value = -9897.000499999998
$ws.Range("N138").Value = -16985.4035

$ws.Range("H141").Value = 5192.3335
$ws.Range("I141").Value = 3161.818
$ws.Range("K141").Value = 9485.454000000002
$ws.Range("M141").Value = -4305.454000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3331.7778
$ws.Range("J2").Value = 4998.5
$ws.Range("L2").Value = 4998.5
$ws.Range("N2").Value = -5224.5

$ws.Range("H45").Value = 3415.6667
$ws.Range("I45").Value = 2595.4
$ws.Range("K45").Value = 2595.4
$ws.Range("M45").Value = -2218.4

$ws.Range("H116").Value = 3331.7778
$ws.Range("J116").Value = 4998.5
$ws.Range("L116").Value = 4998.5
$ws.Range("N116").Value = -9586.5

$ws.Range("H122").Value = 2202.353
$ws.Range("I122").Value = 1962.7333
$ws.Range("K122").Value = 5888.199900000001
$ws.Range("M122").Value = -3438.199900000001

$ws.Range("H132").Value = 2228.75
$ws.Range("I132").Value = 2277
$ws.Range("K132").Value = 6831
$ws.Range("M132").Value = -4301

$ws.Range("H135").Value = 81632.664
$ws.Range("J135").Value = 81632.664
$ws.Range("L135").Value = 81632.664
$ws.Range("N135").Value = -91772.664

$ws.Range("H138").Value = 66590
$ws.Range("J138").Value = 66590
$ws.Range("L138").Value = 66590
$ws.Range("N138").Value = -76870

$ws.Range("H141").Value = 65815
$ws.Range("J141").Value = 65815
$ws.Range("L141").Value = 65815
$ws.Range("N141").Value = -76175

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3331.7778
$ws.Range("J3").Value = 4998.5
$ws.Range("L3").Value = 4998.5
$ws.Range("N3").Value = -5226.5

$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 53583.5
$ws.Range("I63").Value = 50246
$ws.Range("J63").Value = 70271
$ws.Range("K63").Value = 50246
$ws.Range("L63").Value = 70271
$ws.Range("M63").Value = -49560
$ws.Range("N63").Value = -71643

$ws.Range("H66").Value = 53583.5
$ws.Range("I66").Value = 50246
$ws.Range("J66").Value = 70271
$ws.Range("K66").Value = 150738
$ws.Range("L66").Value = 210813
$ws.Range("M66").Value = -147306
$ws.Range("N66").Value = -217677

$ws.Range("H105").Value = 2902.3635
$ws.Range("I105").Value = 1777.8
$ws.Range("J105").Value = 3839.5
$ws.Range("K105").Value = 1777.8
$ws.Range("L105").Value = 3839.5
$ws.Range("M105").Value = -30.79999999999995
$ws.Range("N105").Value = -7333.5

$ws.Range("H108").Value = 83333.336
$ws.Range("J108").Value = 83333.336
$ws.Range("L108").Value = 83333.336
$ws.Range("N108").Value = -91013.336

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H111").Value = 59999
$ws.Range("J111").Value = 59999
$ws.Range("L111").Value = 59999
$ws.Range("N111").Value = -68179

$ws.Range("H134").Value = 3210.4285
$ws.Range("I134").Value = 3226.6155
$ws.Range("K134").Value = 9679.8465
$ws.Range("M134").Value = -7144.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 22286.715
$ws.Range("I62").Value = 19202.2
$ws.Range("J62").Value = 29998
$ws.Range("K62").Value = 57606.60000000001
$ws.Range("L62").Value = 89994
$ws.Range("M62").Value = -56920.60000000001
$ws.Range("N62").Value = -91366

$ws.Range("H65").Value = 22286.715
$ws.Range("I65").Value = 19202.2
$ws.Range("J65").Value = 29998
$ws.Range("K65").Value = 172819.8
$ws.Range("L65").Value = 269982
$ws.Range("M65").Value = -169387.8
$ws.Range("N65").Value = -276846

$ws.Range("H98").Value = 659
$ws.Range("I98").Value = 673.75
$ws.Range("K98").Value = 2021.25
$ws.Range("M98").Value = -523.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 41000
$ws.Range("I58").Value = 41000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 41000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -40723
$ws.Range("N58").ClearContents()

$ws.Range("H64").Value = 55824.57
$ws.Range("I64").Value = 50246
$ws.Range("J64").Value = 69771
$ws.Range("K64").Value = 50246
$ws.Range("L64").Value = 69771
$ws.Range("M64").Value = -49998
$ws.Range("N64").Value = -70267

$ws.Range("H67").Value = 55824.57
$ws.Range("I67").Value = 50246
$ws.Range("J67").Value = 69771
$ws.Range("K67").Value = 50246
$ws.Range("L67").Value = 69771
$ws.Range("M67").Value = -49388
$ws.Range("N67").Value = -71487

$ws.Range("H101").Value = 39367
$ws.Range("J101").Value = 39367
$ws.Range("L101").Value = 39367
$ws.Range("N101").Value = -45857

$ws.Range("H126").Value = 2669.3333
$ws.Range("J126").Value = 2998
$ws.Range("L126").Value = 8994
$ws.Range("N126").Value = -13934

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -860

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H40").Value = 3353.7778
$ws.Range("I40").Value = 3438.6
$ws.Range("K40").Value = 3438.6
$ws.Range("M40").Value = -3302.6

$ws.Range("H100").Value = 2833.25
$ws.Range("I100").Value = 1999.5
$ws.Range("K100").Value = 1999.5
$ws.Range("M100").Value = -1458.5

$ws.Range("H131").Value = 86499.664
$ws.Range("J131").Value = 86499.664
$ws.Range("L131").Value = 86499.664
$ws.Range("N131").Value = -96579.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H127").Value = 48750
$ws.Range("I127").Value = 30000
$ws.Range("K127").Value = 30000
$ws.Range("M127").Value = -25040
